# Regenerate the "K" column (column G) values for rows 2-70 on Sheet1.
# This mirrors the upstream data-regeneration commit that recalculated
# the K (strike count) values and rewrote save_data with the new s_vals.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row (2-70) -> new value for column G ("K")
$newK = @{
    2  = 1
    3  = 2
    4  = 0
    5  = 1
    6  = 2
    7  = 0
    8  = 1
    9  = 0
    10 = 0
    11 = 0
    12 = 1
    13 = 1
    14 = 0
    15 = 0
    16 = 1
    17 = 0
    18 = 1
    19 = 0
    20 = 1
    21 = 0
    22 = 2
    23 = 1
    24 = 0
    25 = 0
    26 = 1
    27 = 0
    28 = 0
    29 = 0
    30 = 0
    31 = 1
    32 = 0
    33 = 0
    34 = 0
    35 = 1
    36 = 0
    37 = 1
    38 = 0
    39 = 1
    40 = 1
    41 = 1
    42 = 1
    43 = 1
    44 = 0
    45 = 1
    46 = 1
    47 = 0
    48 = 2
    49 = 0
    50 = 2
    51 = 3
    52 = 1
    53 = 0
    54 = 2
    55 = 1
    56 = 1
    57 = 1
    58 = 0
    59 = 1
    60 = 0
    61 = 0
    62 = 0
    63 = 2
    64 = 1
    65 = 0
    66 = 0
    67 = 1
    68 = 1
    69 = 0
    70 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
